# Update answer table values in place using Find & Replace.
$d = $word.ActiveDocument

$replacements = @(
    @("926×3=2778", "534×3=1602"),
    @("175×5=875",  "234×2=468"),
    @("541×4=2164", "648×9=5832"),
    @("766×8=6128", "169×7=1183"),
    @("948×8=7584", "604×7=4228"),
    @("554×4=2216", "889×9=8001"),
    @("897×7=6279", "696×7=4872"),
    @("812×4=3248", "193×4=772"),
    @("118×6=708",  "687×7=4809"),
    @("767×3=2301", "800×9=7200"),
    @("852×5=4260", "974×6=5844"),
    @("904×8=7232", "508×2=1016"),
    @("247×7=1729", "965×4=3860"),
    @("203×8=1624", "376×9=3384"),
    @("742×4=2968", "705×8=5640"),
    @("995×9=8955", "105×9=945"),
    @("707×3=2121", "195×8=1560"),
    @("405×5=2025", "901×2=1802"),
    @("146×6=876",  "289×6=1734"),
    @("417×5=2085", "966×6=5796"),
    @("760×7=5320", "780×8=6240"),
    @("650×4=2600", "399×9=3591"),
    @("864×3=2592", "227×2=454"),
    @("994×4=3976", "284×4=1136"),
    @("693×7=4851", "549×3=1647")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
